# Auto-applies F-column ("想去人数" / interest-count) updates
# scraped-data refresh, matching commit "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 444
$ws.Range("F3").Value = 2773
$ws.Range("F4").Value = 1327
$ws.Range("F8").Value = 44
$ws.Range("F9").Value = 601
$ws.Range("F10").Value = 278
$ws.Range("F12").Value = 11538
$ws.Range("F13").Value = 6565
$ws.Range("F15").Value = 16
$ws.Range("F16").Value = 414
$ws.Range("F17").Value = 255
$ws.Range("F21").Value = 62
$ws.Range("F22").Value = 263
$ws.Range("F23").Value = 920
$ws.Range("F24").Value = 3634
$ws.Range("F25").Value = 55
$ws.Range("F28").Value = 162
$ws.Range("F30").Value = 17
$ws.Range("F31").Value = 266
$ws.Range("F32").Value = 294
$ws.Range("F33").Value = 4996
$ws.Range("F35").Value = 1229
$ws.Range("F36").Value = 226
$ws.Range("F37").Value = 422
$ws.Range("F39").Value = 531

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9026
$ws.Range("F3").Value = 494
$ws.Range("F4").Value = 1814

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9026
$ws.Range("F3").Value = 494
$ws.Range("F4").Value = 1814
$ws.Range("F5").Value = 444
$ws.Range("F6").Value = 2773
$ws.Range("F10").Value = 1327
$ws.Range("F14").Value = 44
$ws.Range("F16").Value = 601
$ws.Range("F17").Value = 278
$ws.Range("F19").Value = 11538
$ws.Range("F23").Value = 16
$ws.Range("F24").Value = 414
$ws.Range("F25").Value = 255
$ws.Range("F29").Value = 62
$ws.Range("F30").Value = 263
$ws.Range("F31").Value = 920
$ws.Range("F32").Value = 3634
$ws.Range("F33").Value = 55
$ws.Range("F35").Value = 162
$ws.Range("F37").Value = 266
$ws.Range("F41").Value = 1229
$ws.Range("F42").Value = 226
$ws.Range("F44").Value = 531
